# aggiunta schermata e logica per i dinamici, da testare
#
# 1) Typography sheet: set the Fallback Character (column F) of row 5 to "%"
# 2) Translation sheet: fill in the new dynamic-screen text rows 339-346
#    (Text ID / Typography Name / Alignment / GB text / Direction)

$wb = $excel.ActiveWorkbook

$wsTypography = $wb.Worksheets.Item("Typography")
$wsTypography.Cells(5, 6).Value = "%"

$wsTranslation = $wb.Worksheets.Item("Translation")

$rows = @(
    @{ Row = 339; Id = "SingleUseId363"; Font = "Medium"; Align = "Left";   Text = "CALIBRATION" },
    @{ Row = 340; Id = "SingleUseId364"; Font = "Medium"; Align = "Left";   Text = "APPS 0%" },
    @{ Row = 341; Id = "SingleUseId365"; Font = "Medium"; Align = "Left";   Text = "APPS 100%" },
    @{ Row = 342; Id = "SingleUseId366"; Font = "Medium"; Align = "Left";   Text = "SW ANGLE" },
    @{ Row = 343; Id = "SingleUseId367"; Font = "Medium"; Align = "Left";   Text = "LINEAR" },
    @{ Row = 344; Id = "SingleUseId368"; Font = "Medium"; Align = "Left";   Text = "LOAD CELL" },
    @{ Row = 345; Id = "SingleUseId369"; Font = "Medium"; Align = "Center"; Text = "<value> DONE" },
    @{ Row = 346; Id = "SingleUseId370"; Font = "Medium"; Align = "Left";   Text = "NOTHING" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $wsTranslation.Cells($rowNum, 2).Value = $r.Id     # B - Text ID
    $wsTranslation.Cells($rowNum, 3).Value = $r.Font    # C - Typography Name
    $wsTranslation.Cells($rowNum, 4).Value = $r.Align   # D - Alignment
    $wsTranslation.Cells($rowNum, 5).Value = $r.Text    # E - GB (text)
    $wsTranslation.Cells($rowNum, 6).Value = "LTR"      # F - Direction
}
